$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("W2").Value = 894
$ws.Range("W3").Value = 99
$ws.Range("W4").Value = 492
$ws.Range("W5").Value = 140
$ws.Range("W6").Value = 589
$ws.Range("W7").Value = 745
$ws.Range("W8").Value = 2959
